$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.721.39"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "3.789.29"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.43"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.07"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "3.788.05"
$ws.Range("E7").Value = "  +1.42%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("E11").Value = "  -1.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("E13").Value = "  -1.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.10"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "4.413.52"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").Value = "3.799.87"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("D18").Value = "67.666.37"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("E21").Value = "  -6.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.44"
$ws.Range("E22").Value = "  -1.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.698"

$ws.Range("E24").Value = "  +5.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.35"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("E27").Value = "  -2.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +4.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.24"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.68"
$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.100"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.995"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.78"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.22"
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.07"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.298"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.63"
$ws.Range("E46").Value = "  +3.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.30"
$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "394.18"
$ws.Range("E48").Value = "  +1.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.76"
$ws.Range("E49").Value = "  +6.86%  "

$ws.Range("E51").Value = "  +7.09%  "
